$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so values such as
# "1.00" or "603.06" are not silently coerced into numbers by Excel,
# matching the inlineStr cells in the source workbook.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '66.010.56'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '3.175.06'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '603.06'
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").Value = '153.89'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.176.86'
$ws.Range("E8").Value = '  -0.61%  '
$ws.Range("D9").Value = '0.544'
$ws.Range("E9").Value = '  +1.72%  '
$ws.Range("D10").Value = '0.158'
$ws.Range("E10").Value = '  -1.67%  '
$ws.Range("D11").Value = '5.63'
$ws.Range("E11").Value = '  -7.77%  '
$ws.Range("D12").Value = '0.505'
$ws.Range("E12").Value = '  -1.74%  '
$ws.Range("D13").Value = '0.0000264'
$ws.Range("E13").Value = '  -2.78%  '
$ws.Range("D14").Value = '38.16'
$ws.Range("E14").Value = '  -2.17%  '
$ws.Range("D15").Value = '3.700.32'
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").Value = '66.083.51'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '7.33'
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("D18").Value = '3.178.34'
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").Value = '0.112'
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("D20").Value = '505.33'
$ws.Range("E20").Value = '  -1.00%  '
$ws.Range("D21").Value = '15.25'
$ws.Range("E21").Value = '  -0.69%  '
$ws.Range("D22").Value = '0.726'
$ws.Range("E22").Value = '  -1.98%  '
$ws.Range("D23").Value = '7.98'
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").Value = '14.73'
$ws.Range("E24").Value = '  -3.26%  '
$ws.Range("D25").Value = '84.21'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").Value = '2.98'
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("D28").Value = '9.09'
$ws.Range("E28").Value = '  -3.32%  '
$ws.Range("D29").Value = '2.37'
$ws.Range("E29").Value = '  +4.58%  '
$ws.Range("D30").Value = '3.02'
$ws.Range("E30").Value = '  +5.10%  '
$ws.Range("D31").Value = '6.98'
$ws.Range("E31").Value = '  +1.38%  '
$ws.Range("D32").Value = '27.81'
$ws.Range("E32").Value = '  -1.65%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.14%  '
$ws.Range("D34").Value = '1.17'
$ws.Range("E34").Value = '  -4.56%  '
$ws.Range("D35").Value = '6.45'
$ws.Range("E35").Value = '  -1.63%  '
$ws.Range("D36").Value = '509.16'
$ws.Range("E36").Value = '  +4.81%  '
$ws.Range("D37").Value = '55.23'
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("D38").Value = '0.0915'
$ws.Range("E38").Value = '  +1.15%  '
$ws.Range("D39").Value = '0.0415'
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("D40").Value = '0.0₃0706'
$ws.Range("E40").Value = '  +8.77%  '
$ws.Range("D41").Value = '0.127'
$ws.Range("E41").Value = '  +4.15%  '
$ws.Range("D42").Value = '8.74'
$ws.Range("E42").Value = '  -1.10%  '
$ws.Range("D43").Value = '2.86'
$ws.Range("E43").Value = '  -2.44%  '
$ws.Range("D44").Value = '0.297'
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("D45").Value = '2.45'
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("D46").Value = '2.828.05'
$ws.Range("E46").Value = '  -3.87%  '
$ws.Range("D47").Value = '27.71'
$ws.Range("E47").Value = '  -2.50%  '
$ws.Range("D48").Value = '0.999'
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("D49").Value = '2.36'
$ws.Range("E49").Value = '  +2.52%  '
$ws.Range("D50").Value = '0.116'
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("D51").Value = '2.65'
$ws.Range("E51").Value = '  +2.66%  '

# Restore the default (unformatted) style so the cells end up identical
# to the original workbook's styling (no explicit style index).
$ws.Range("D2:E51").Style = "Normal"
